$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2

$d.Content.Find.Execute(
    "{molecular dynamics (MD)|simulation} simulations were performed using",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "{molecular dynamics (MD)|simulation} were performed using", 2)

$d.Content.Find.Execute(
    "{ff14SB|force field} force field was used.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "{ff14SB|force field} was used.", 2)

$d.Content.Find.Execute(
    "{ff19SB|force field} force field was used.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "{ff19SB|force field} was used.", 2)

$d.Content.Find.Execute(
    "{LIPID14|force field} force field",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "{LIPID14|force field}.", 2)

$d.Content.Find.Execute(
    "{NVT|MD} MD simulations were performed.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "{NVT|MD} simulations were performed.", 2)
